$d = $word.ActiveDocument

$replacements = @(
    @("2025-02-22 Saturday", "2025-02-23 Sunday"),
    @("11×96=1056", "15×19=285"),
    @("84×66=5544", "50×79=3950"),
    @("76×36=2736", "57×75=4275"),
    @("51×41=2091", "91×44=4004"),
    @("17×92=1564", "56×24=1344"),
    @("75×59=4425", "86×73=6278"),
    @("15×73=1095", "39×58=2262"),
    @("58×45=2610", "71×77=5467"),
    @("49×53=2597", "55×34=1870"),
    @("24×97=2328", "69×65=4485"),
    @("44×81=3564", "16×88=1408"),
    @("46×48=2208", "22×21=462"),
    @("14×19=266", "32×98=3136"),
    @("89×38=3382", "48×68=3264"),
    @("15×68=1020", "66×15=990"),
    @("39×90=3510", "98×16=1568"),
    @("40×60=2400", "94×61=5734"),
    @("43×32=1376", "66×72=4752"),
    @("25×67=1675", "68×27=1836"),
    @("99×97=9603", "45×15=675"),
    @("20×79=1580", "69×99=6831"),
    @("49×93=4557", "51×15=765"),
    @("53×17=901", "35×95=3325"),
    @("72×18=1296", "54×13=702"),
    @("34×67=2278", "58×49=2842")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
